$wb = $excel.ActiveWorkbook

# --- Sheet1: selection moves from C2 to E11 (tabSelected cleared once Sheet4 is activated later) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E11").Select()

# --- Sheet4: C2 keeps its value "(build 39769)"; selection moves from D2 to C2 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("C2").Value = "(build 39769)"
$ws4.Range("C2").Select()

# --- New Sheet5, appended after Sheet4 ---
$ws5 = $wb.Worksheets.Add($null, $ws4)

# Column headers (row 1), filled green
$ws5.Range("A1").Value = "URL"
$ws5.Range("B1").Value = "Username"
$ws5.Range("C1").Value = "Password"
$ws5.Range("D1").Value = "FirstName"
$ws5.Range("E1").Value = "MiddleName"
$ws5.Range("F1").Value = "LastName"
$ws5.Range("G1").Value = "email"
$ws5.Range("H1").Value = "UserName"
$ws5.Range("I1").Value = "Password"
$ws5.Range("J1").Value = "Repassword"
$ws5.Range("K1").Value = "Working schedule"
$ws5.Range("L1").Value = "Search Name"
$ws5.Range("M1").Value = "title"
$ws5.Range("A1:M1").Interior.Color = 5296274

# Data row 2
$ws5.Range("A2").Value = "http://localhost/login.do"
$ws5.Range("B2").Value = "admin"
$ws5.Range("C2").Value = "manager"
$ws5.Range("D2").Value = "shekhar"
$ws5.Range("E2").Value = "m"
$ws5.Range("F2").Value = "pakale"
$ws5.Range("G2").Value = "shekhar.pakale@gmail.com"
$ws5.Hyperlinks.Add($ws5.Range("G2"), "mailto:shekhar.pakale@gmail.com")
$ws5.Range("H2").Value = "shekhar123"
$ws5.Range("I2").Value = "pass123"
$ws5.Range("J2").Value = "pass123"
$ws5.Range("K2").Value = 10
$ws5.Range("K2").HorizontalAlignment = -4131
$ws5.Range("K2").VerticalAlignment = -4160
$ws5.Range("L2").Value = "shekhar"
$ws5.Range("M2").Value = "actiTIME - User List"

# Data row 3 (partial)
$ws5.Range("B3").Value = "shekharp123"
$ws5.Range("C3").Value = "pass123"

# Column widths (best-fit approximation)
$ws5.Columns("A:M").AutoFit()

$ws5.Range("B4").Select()

# --- Activate Sheet4 as the final selected sheet/tab ---
$ws4.Activate()
